$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = '66.690.47'
$ws.Range("E2").Value = '  +5.37%  '

# Row 3
$ws.Range("D3").Value = '3.719.03'
$ws.Range("E3").Value = '  +7.40%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
Set-TextValue ($ws.Range("D5")) '423.08'
$ws.Range("E5").Value = '  +1.89%  '

# Row 6
Set-TextValue ($ws.Range("D6")) '131.75'
$ws.Range("E6").Value = '  +1.85%  '

# Row 7
$ws.Range("D7").Value = '3.714.34'
$ws.Range("E7").Value = '  +7.42%  '

# Row 8
Set-TextValue ($ws.Range("D8")) '0.639'
$ws.Range("E8").Value = '  +1.73%  '

# Row 9
$ws.Range("E9").Value = '  +0.04%  '

# Row 10
Set-TextValue ($ws.Range("D10")) '0.761'
$ws.Range("E10").Value = '  -0.76%  '

# Row 11
Set-TextValue ($ws.Range("D11")) '0.179'
$ws.Range("E11").Value = '  +14.74%  '

# Row 12
Set-TextValue ($ws.Range("D12")) '0.0000378'
$ws.Range("E12").Value = '  +64.70%  '

# Row 13
Set-TextValue ($ws.Range("D13")) '42.60'
$ws.Range("E13").Value = '  +1.09%  '

# Row 14
Set-TextValue ($ws.Range("D14")) '10.24'
$ws.Range("E14").Value = '  +4.18%  '

# Row 15
$ws.Range("D15").Value = '4.299.11'
$ws.Range("E15").Value = '  +6.76%  '

# Row 16
$ws.Range("E16").Value = '  -0.05%  '

# Row 17
Set-TextValue ($ws.Range("D17")) '20.80'
$ws.Range("E17").Value = '  +2.68%  '

# Row 18
$ws.Range("D18").Value = '3.670.62'
$ws.Range("E18").Value = '  +5.46%  '

# Row 19
Set-TextValue ($ws.Range("D19")) '12.86'
$ws.Range("E19").Value = '  +3.66%  '

# Row 20
Set-TextValue ($ws.Range("D20")) '1.12'
$ws.Range("E20").Value = '  +2.63%  '

# Row 21
$ws.Range("D21").Value = '66.837.77'
$ws.Range("E21").Value = '  +5.57%  '

# Row 22
Set-TextValue ($ws.Range("D22")) '444.07'
$ws.Range("E22").Value = '  -2.97%  '

# Row 23
Set-TextValue ($ws.Range("D23")) '15.42'
$ws.Range("E23").Value = '  +15.07%  '

# Row 24
Set-TextValue ($ws.Range("D24")) '89.24'
$ws.Range("E24").Value = '  -0.95%  '

# Row 25
$ws.Range("E25").Value = '  -4.35%  '

# Row 26
Set-TextValue ($ws.Range("D26")) '37.85'
$ws.Range("E26").Value = '  +13.37%  '

# Row 27
Set-TextValue ($ws.Range("D27")) '10.16'
$ws.Range("E27").Value = '  +0.07%  '

# Row 28
$ws.Range("E28").Value = '  +0.13%  '

# Row 29
$ws.Range("E29").Value = '  +5.19%  '

# Row 30
Set-TextValue ($ws.Range("D30")) '12.54'
$ws.Range("E30").Value = '  +1.09%  '

# Row 31
Set-TextValue ($ws.Range("D31")) '2.77'
$ws.Range("E31").Value = '  +3.75%  '

# Row 32
Set-TextValue ($ws.Range("D32")) '0.121'
$ws.Range("E32").Value = '  +7.61%  '

# Row 33
Set-TextValue ($ws.Range("D33")) '7.24'
$ws.Range("E33").Value = '  -3.47%  '

# Row 34
Set-TextValue ($ws.Range("D34")) '41.93'
$ws.Range("E34").Value = '  +5.17%  '

# Row 35
Set-TextValue ($ws.Range("D35")) '0.163'
$ws.Range("E35").Value = '  -1.32%  '

# Row 36
$ws.Range("E36").Value = '  +0.37%  '

# Row 37
Set-TextValue ($ws.Range("D37")) '56.40'
$ws.Range("E37").Value = '  -2.04%  '

# Row 38
Set-TextValue ($ws.Range("D38")) '0.0488'
$ws.Range("E38").Value = '  +0.18%  '

# Row 39
$ws.Range("D39").Value = '0.0₃0738'
$ws.Range("E39").Value = '  +16.20%  '

# Row 40
Set-TextValue ($ws.Range("D40")) '3.08'
$ws.Range("E40").Value = '  +32.10%  '

# Row 41
Set-TextValue ($ws.Range("D41")) '0.145'
$ws.Range("E41").Value = '  +5.64%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue ($ws.Range("D42")) '28.48'
$ws.Range("E42").Value = '  +29.60%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue ($ws.Range("D43")) '0.997'
$ws.Range("E43").Value = '  -0.22%  '

# Row 44
$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue ($ws.Range("D44")) '3.47'
$ws.Range("E44").Value = '  +4.00%  '

# Row 46
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue ($ws.Range("D46")) '145.36'
$ws.Range("E46").Value = '  -0.31%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue ($ws.Range("D47")) '2.91'
$ws.Range("E47").Value = '  -5.44%  '

# Row 48
Set-TextValue ($ws.Range("D48")) '2.65'
$ws.Range("E48").Value = '  -6.05%  '

# Row 49
Set-TextValue ($ws.Range("D49")) '4.35'
$ws.Range("E49").Value = '  -1.92%  '

# Row 50
$ws.Range("E50").Value = '  -4.25%  '

# Row 51
Set-TextValue ($ws.Range("D51")) '0.159'
$ws.Range("E51").Value = '  +14.77%  '

